# "Generate Report for Handoff"
# Updates the localization-status report: the Status moves from
# "In Translation" to "Ready for handoff", the corresponding handoff
# timestamps are refreshed, and the "Status" columns are widened to fit
# the new (longer) text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" -----------------
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value     = "Ready for handoff"
$wsDeDe.Range("C2").Value     = "Ready for handoff"

# --- Refreshed handoff datetimes --------------------------------------
# Overview "Latest HO Xliff Generate Date" and de-de "Latest Handoff Datetime"
$wsOverview.Range("G2").Value = "2016-09-03 08:42:40"
$wsDeDe.Range("H2").Value     = "2016-09-03 08:42:40"

# zh-cn "Latest Handoff Datetime"
$wsZhCn.Range("H2").Value = "2016-09-03 08:42:35"

# --- Widen the "Status" columns to fit "Ready for handoff" ------------
# Target stored width is ~17.216; the ColumnWidth property snaps to a
# 1/6-character grid, so 16.3333 is the closest input that lands on the
# nearest achievable stored width.
$wsOverview.Columns.Item(5).ColumnWidth = 16.3333333333333
$wsOverview.Columns.Item(6).ColumnWidth = 16.3333333333333
$wsZhCn.Columns.Item(3).ColumnWidth     = 16.3333333333333
$wsDeDe.Columns.Item(3).ColumnWidth     = 16.3333333333333
